$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "sd_Y"

$values = @(
    0,
    0.102617110414406,
    0.0904249124830942,
    0.0948773861746983,
    0,
    0.131305920744823,
    0.0782046928103035,
    0.0389651517486449,
    0,
    0.193429739786635,
    0.203129339660775,
    0.0339220782532088,
    0.031600616644903,
    0.172502458799703,
    0.141525516201864,
    0
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
